$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 176 (shifts existing rows 176-222 down to 177-223)
$ws.Rows("176:176").Insert()

# Populate the newly inserted row 176 with the new weekly price observation
$ws.Range("A176").Value = 4
$ws.Range("B176").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C176").Value = "Los Lagos"
$ws.Range("D176").Value = 44642
$ws.Range("E176").Value = 10
$ws.Range("F176").Value = 100112017
$ws.Range("G176").Value = "Apio"
$ws.Range("H176").Value = "Americana (o)"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 45
$ws.Range("K176").Value = 12000
$ws.Range("L176").Value = 12000
$ws.Range("M176").Value = 12000
$ws.Range("N176").Value = "$/docena de matas"
$ws.Range("O176").Value = "Región de Coquimbo"
$ws.Range("P176").Value = 2000
$ws.Range("Q176").Value = 6
$ws.Range("R176").Value = "Hortaliza"
